# Apply updates to the absence rate workbook:
#  1. Rename the elementary schools in column B, stripping the trailing
#     " ES" suffix (leaving two trailing spaces in its place).
#  2. Update the rounded "Absence Rate" values in column D.
#  3. Select cell D3 to match the saved selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename school names (column B) across all data rows.
$nameMap = @{
    "Forest Grove ES"  = "Forest Grove  "
    "Guilford ES"      = "Guilford  "
    "Rolling Ridge ES" = "Rolling Ridge  "
    "Sterling ES"      = "Sterling  "
    "Sugarland ES"     = "Sugarland  "
    "Sully ES"         = "Sully  "
}

for ($row = 2; $row -le 25; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $current = $cell.Value2
    if ($nameMap.ContainsKey($current)) {
        $cell.Value = $nameMap[$current]
    }
}

# 2. Update the Absence Rate values (column D) with new rounded figures.
$absenceRates = @{
    2  = 0.041
    3  = 0.044
    4  = 0.052
    5  = 0.059
    6  = 0.061
    7  = 0.067
    8  = 0.056
    9  = 0.058
    10 = 0.053
    11 = 0.069
    12 = 0.066
    13 = 0.071
    14 = 0.045
    15 = 0.049
    16 = 0.051
    17 = 0.053
    18 = 0.065
    19 = 0.078
    20 = 0.063
    21 = 0.064
    22 = 0.053
    23 = 0.067
    24 = 0.076
    25 = 0.077
}

foreach ($row in $absenceRates.Keys) {
    $ws.Cells.Item($row, 4).Value = $absenceRates[$row]
}

# 3. Match the active selection recorded in the saved file.
$ws.Range("D3").Select()
